$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.422.03'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '3.924.92'
$ws.Range('E3').Value = '  +3.99%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'470.98"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.80%  '
$ws.Range('D6').Value = "'145.19"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.34%  '
$ws.Range('D7').Value = "'0.625"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = "'0.998"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.733"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').Value = "'0.167"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.50%  '
$ws.Range('D11').Value = "'0.0000342"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.60%  '
$ws.Range('D12').Value = "'43.36"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '4.547.49'
$ws.Range('E13').Value = '  +4.15%  '
$ws.Range('D14').Value = "'10.44"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = "'15.04"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '3.917.58'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = "'19.89"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  +3.60%  '
$ws.Range('D20').Value = '67.631.75'
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').Value = "'434.65"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.57%  '
$ws.Range('D22').Value = "'14.65"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('D23').Value = "'3.35"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('D24').Value = "'88.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.17%  '
$ws.Range('D25').Value = "'38.75"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.18%  '
$ws.Range('D26').Value = "'3.54"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.72%  '
$ws.Range('E27').Value = '  +4.00%  '
$ws.Range('D28').Value = "'10.14"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('D29').Value = "'9.64"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.53%  '
$ws.Range('D30').Value = "'723.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('D31').Value = "'13.64"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = "'0.132"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('D34').Value = "'43.03"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.83%  '
$ws.Range('E35').Value = '  +4.27%  '
$ws.Range('D36').Value = "'57.76"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('D37').Value = '0.0₃0806'
$ws.Range('E37').Value = '  +21.07%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  -6.92%  '
$ws.Range('D40').Value = "'0.0475"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').Value = "'3.06"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = "'2.58"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.43%  '
$ws.Range('E44').Value = '  +2.17%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'1.00"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('E47').Value = '  +5.99%  '
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').Value = "'3.16"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').Value = "'145.11"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.47%  '
$ws.Range('E51').Value = '  +3.83%  '
